# Add a new row 5 to the meta-sheet for release/6.0.3,
# matching the existing rows' pattern (A = release name, B:D = "X").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "release/6.0.3"
$ws.Range("B5").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("D5").Value = "X"
